$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.958.50'
$ws.Range('E2').Value = '  +3.62%  '
$ws.Range('D3').Value = '2.422.86'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'315.61"
$ws.Range('E5').Value = '  +2.90%  '
$ws.Range('D6').Value = "'102.63"
$ws.Range('E6').Value = '  +5.56%  '
$ws.Range('D7').Value = "'0.513"
$ws.Range('E7').Value = '  +1.13%  '
$ws.Range('D8').Value = "'0.999"
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = "'0.524"
$ws.Range('E9').Value = '  +7.25%  '
$ws.Range('D10').Value = "'35.38"
$ws.Range('E10').Value = '  +1.18%  '
$ws.Range('D11').Value = "'0.0801"
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('E12').Value = '  -2.22%  '
$ws.Range('D13').Value = "'18.22"
$ws.Range('E13').Value = '  -1.90%  '
$ws.Range('D14').Value = "'6.98"
$ws.Range('E14').Value = '  +1.38%  '
$ws.Range('D15').Value = '2.801.26'
$ws.Range('E15').Value = '  +0.85%  '
$ws.Range('D16').Value = '2.431.02'
$ws.Range('E16').Value = '  +1.91%  '
$ws.Range('E17').Value = '  +0.67%  '
$ws.Range('D18').Value = '44.872.54'
$ws.Range('D19').Value = "'12.22"
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('D20').Value = "'6.36"
$ws.Range('E20').Value = '  -0.91%  '
$ws.Range('D21').Value = '0.0₃0922'
$ws.Range('E21').Value = '  +2.12%  '
$ws.Range('D22').Value = "'68.68"
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').Value = "'243.44"
$ws.Range('E23').Value = '  +2.42%  '
$ws.Range('D24').Value = "'2.26"
$ws.Range('E24').Value = '  +0.48%  '
$ws.Range('E25').Value = '  +0.93%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').Value = "'25.18"
$ws.Range('E27').Value = '  +0.79%  '
$ws.Range('D28').Value = "'2.27"
$ws.Range('E28').Value = '  -2.78%  '
$ws.Range('D29').Value = "'9.54"
$ws.Range('E29').Value = '  +1.34%  '
$ws.Range('D30').Value = "'49.07"
$ws.Range('E30').Value = '  +1.74%  '
$ws.Range('D31').Value = "'32.73"
$ws.Range('E31').Value = '  +1.52%  '
$ws.Range('D32').Value = "'19.83"
$ws.Range('E32').Value = '  +7.35%  '
$ws.Range('D33').Value = "'0.124"
$ws.Range('E33').Value = '  +7.68%  '
$ws.Range('D34').Value = "'5.22"
$ws.Range('E34').Value = '  +1.42%  '
$ws.Range('E35').Value = '  +0.41%  '
$ws.Range('D36').Value = "'0.0761"
$ws.Range('E36').Value = '  +1.75%  '
$ws.Range('D37').Value = "'1.87"
$ws.Range('E37').Value = '  -1.23%  '
$ws.Range('E38').Value = '  +0.39%  '
$ws.Range('D39').Value = "'2.86"
$ws.Range('E39').Value = '  -6.01%  '
$ws.Range('D40').Value = "'122.22"
$ws.Range('E40').Value = '  -6.81%  '
$ws.Range('D41').Value = "'2.21"
$ws.Range('E41').Value = '  -2.96%  '
$ws.Range('E42').Value = '  +1.15%  '
$ws.Range('D43').Value = "'20.72"
$ws.Range('E43').Value = '  -2.88%  '
$ws.Range('E44').Value = '  +2.23%  '
$ws.Range('D45').Value = '1.932.11'
$ws.Range('E45').Value = '  -0.93%  '
$ws.Range('D46').Value = "'2.10"
$ws.Range('E46').Value = '  -2.85%  '
$ws.Range('D47').Value = "'2.92"
$ws.Range('E47').Value = '  +4.34%  '
$ws.Range('D48').Value = "'9.18"
$ws.Range('E48').Value = '  -2.74%  '
$ws.Range('D49').Value = "'1.78"
$ws.Range('E49').Value = '  +15.20%  '
$ws.Range('D50').Value = "'76.33"
$ws.Range('E50').Value = '  +5.81%  '
$ws.Range('D51').Value = "'53.91"
$ws.Range('E51').Value = '  +2.70%  '
